# Update crypto price/volume data per Fri Jun 21 15:40:22 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.794.26"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.496.67"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.56"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.54"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.486.48"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.13"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.095.25"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.42"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.498.18"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.963.03"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.65"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.95"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.577"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.643.64"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.96"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.512.60"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.51"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.94"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.66"
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0796"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.45"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.807"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.22"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.29"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.420.27"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.890"
$ws.Range("E51").Value = "  -0.86%  "
